$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 736.4286
$ws.Range("I19").Value = 880
$ws.Range("J19").Value = 592.8570999999999
$ws.Range("K19").Value = 880
$ws.Range("L19").Value = 592.8570999999999
$ws.Range("M19").Value = -705
$ws.Range("N19").Value = -942.8570999999999
$ws.Range("H51").Value = 7844.1113
$ws.Range("J51").Value = 8083.1665
$ws.Range("L51").Value = 8083.1665
$ws.Range("N51").Value = -9051.166499999999
$ws.Range("H62").Value = 4116.346
$ws.Range("I62").Value = 3091.6667
$ws.Range("K62").Value = 3091.6667
$ws.Range("M62").Value = -2467.6667
$ws.Range("H65").Value = 4116.346
$ws.Range("I65").Value = 3091.6667
$ws.Range("K65").Value = 15458.3335
$ws.Range("M65").Value = -12338.3335
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()
$ws.Range("H112").Value = 1625.1538
$ws.Range("J112").Value = 3000
$ws.Range("L112").Value = 9000
$ws.Range("N112").Value = -11216
$ws.Range("H138").Value = 2618.9556
$ws.Range("I138").Value = 1960.0454
$ws.Range("K138").Value = 5880.1362
$ws.Range("M138").Value = -740.1361999999999
$ws.Range("H141").Value = 7199.9
$ws.Range("I141").Value = 1999
$ws.Range("K141").Value = 5997
$ws.Range("M141").Value = -817

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 55559880
$ws.Range("I45").Value = 111112180
$ws.Range("J45").Value = 7584.4443
$ws.Range("K45").Value = 111112180
$ws.Range("L45").Value = 7584.4443
$ws.Range("M45").Value = -111111803
$ws.Range("N45").Value = -8338.444299999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2401
$ws.Range("I20").Value = 1957.5
$ws.Range("J20").Value = 3288
$ws.Range("K20").Value = 1957.5
$ws.Range("L20").Value = 3288
$ws.Range("M20").Value = -1710.5
$ws.Range("N20").Value = -3782
$ws.Range("H86").Value = 5368.643
$ws.Range("I86").Value = 3020.5
$ws.Range("J86").Value = 11239
$ws.Range("K86").Value = 3020.5
$ws.Range("L86").Value = 11239
$ws.Range("M86").Value = -1897.5
$ws.Range("N86").Value = -13485
$ws.Range("H89").Value = 5368.643
$ws.Range("I89").Value = 3020.5
$ws.Range("J89").Value = 11239
$ws.Range("K89").Value = 15102.5
$ws.Range("L89").Value = 56195
$ws.Range("M89").Value = -9486.5
$ws.Range("N89").Value = -67427
$ws.Range("H107").Value = 1300.5834
$ws.Range("I107").Value = 1360.8
$ws.Range("K107").Value = 1360.8
$ws.Range("M107").Value = 559.2
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()
$ws.Range("H140").Value = 69776.8
$ws.Range("J140").Value = 69776.8
$ws.Range("L140").Value = 69776.8
$ws.Range("N140").Value = -80136.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 57360.75
$ws.Range("I31").Value = 6354
$ws.Range("K31").Value = 6354
$ws.Range("M31").Value = -6059
$ws.Range("H34").Value = 57360.75
$ws.Range("I34").Value = 6354
$ws.Range("K34").Value = 6354
$ws.Range("M34").Value = -6152
$ws.Range("H132").Value = 8615.147999999999
$ws.Range("I132").Value = 8626.653
$ws.Range("J132").Value = 8502.4
$ws.Range("K132").Value = 25879.959
$ws.Range("L132").Value = 25507.2
$ws.Range("M132").Value = -23349.959
$ws.Range("N132").Value = -30567.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 67052.44500000001
$ws.Range("J2").Value = 67052.44500000001
$ws.Range("L2").Value = 402314.67
$ws.Range("N2").Value = -402540.67
$ws.Range("H23").Value = 371.81818
$ws.Range("I23").Value = 100
$ws.Range("J23").Value = 527.1429000000001
$ws.Range("K23").Value = 300
$ws.Range("L23").Value = 1581.4287
$ws.Range("M23").Value = -65
$ws.Range("N23").Value = -2051.4287
$ws.Range("H34").Value = 3211.3076
$ws.Range("J34").Value = 5133.3335
$ws.Range("L34").Value = 15400.0005
$ws.Range("N34").Value = -15568.0005
$ws.Range("H137").Value = 4544.8
$ws.Range("J137").Value = 6310.5
$ws.Range("L137").Value = 18931.5
$ws.Range("N137").Value = -29131.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 14496.417
$ws.Range("I70").Value = 5536
$ws.Range("J70").Value = 36257.43
$ws.Range("K70").Value = 5536
$ws.Range("L70").Value = 36257.43
$ws.Range("M70").Value = -5266
$ws.Range("N70").Value = -36797.43
$ws.Range("H73").Value = 14496.417
$ws.Range("I73").Value = 5536
$ws.Range("J73").Value = 36257.43
$ws.Range("K73").Value = 5536
$ws.Range("L73").Value = 36257.43
$ws.Range("M73").Value = -4600
$ws.Range("N73").Value = -38129.43
$ws.Range("H102").Value = 3242.2
$ws.Range("I102").Value = 3437.3333
$ws.Range("K102").Value = 3437.3333
$ws.Range("M102").Value = -1815.3333
$ws.Range("H126").Value = 3675.0435
$ws.Range("I126").Value = 2926.9473
$ws.Range("K126").Value = 8780.841899999999
$ws.Range("M126").Value = -6310.841899999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1702.5
$ws.Range("I16").Value = 1448.1818
$ws.Range("K16").Value = 1448.1818
$ws.Range("M16").Value = -1278.1818
$ws.Range("H22").Value = 4910.8335
$ws.Range("I22").Value = 2492
$ws.Range("J22").Value = 8297.200000000001
$ws.Range("K22").Value = 2492
$ws.Range("L22").Value = 8297.200000000001
$ws.Range("M22").Value = -2197
$ws.Range("N22").Value = -8887.200000000001
$ws.Range("H26").Value = 18336.666
$ws.Range("I26").Value = 10000
$ws.Range("J26").Value = 22505
$ws.Range("K26").Value = 10000
$ws.Range("L26").Value = 22505
$ws.Range("M26").Value = -9705
$ws.Range("N26").Value = -23095
$ws.Range("H27").Value = 4910.8335
$ws.Range("I27").Value = 2492
$ws.Range("J27").Value = 8297.200000000001
$ws.Range("K27").Value = 2492
$ws.Range("L27").Value = 8297.200000000001
$ws.Range("M27").Value = -2385
$ws.Range("N27").Value = -8511.200000000001
$ws.Range("H40").Value = 6403.7744
$ws.Range("I40").Value = 6413.609
$ws.Range("K40").Value = 6413.609
$ws.Range("M40").Value = -6277.609
$ws.Range("H55").Value = 2502442.2
$ws.Range("J55").Value = 3947.3
$ws.Range("L55").Value = 3947.3
$ws.Range("N55").Value = -4293.3
$ws.Range("H82").Value = 4861.4443
$ws.Range("I82").Value = 3375.4707
$ws.Range("K82").Value = 3375.4707
$ws.Range("M82").Value = -3014.4707
$ws.Range("H85").Value = 4861.4443
$ws.Range("I85").Value = 3375.4707
$ws.Range("K85").Value = 3375.4707
$ws.Range("M85").Value = -2127.4707

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1024.1111
$ws.Range("J113").Value = 1064.875
$ws.Range("L113").Value = 3194.625
$ws.Range("N113").Value = -7534.625
$ws.Range("H132").Value = 5453.7036
$ws.Range("I132").Value = 5453.1875
$ws.Range("K132").Value = 16359.5625
$ws.Range("M132").Value = -13829.5625
